{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// --- 1) Collapse the \"SwiftPay\"-split runs (removes spell-check\n//        proofErr markup by rewriting each paragraph's text as one run).\n\n// Title paragraph: \"Phase 3: Software Architecture Design \u2013 SwiftPay\"\nconst titlePara = paragraphs.items.find(\n  (p) => p.text === \"Phase 3: Software Architecture Design \u2013 SwiftPay\"\n);\nif (titlePara) {\n  titlePara.clear();\n  titlePara.insertText(\n    \"Phase 3: Software Architecture Design \u2013 SwiftPay\",\n    Word.InsertLocation.start\n  );\n}\n\n// Overview paragraph.\nconst overviewPara = paragraphs.items.find(\n  (p) => p.text.indexOf(\"The system architecture of\") === 0\n);\nif (overviewPara) {\n  overviewPara.clear();\n  overviewPara.insertText(\n    \"The system architecture of SwiftPay outlines how different components of the mobile application will interact to deliver financial services such as registration, login, balance inquiry, money transfer, and more. The system is designed to be scalable, secure, and maintainable.\",\n    Word.InsertLocation.start\n  );\n}\n\n// \"SwiftPay follows a three-tier architecture, which includes:\" paragraph.\nconst tierPara = paragraphs.items.find(\n  (p) => p.text.indexOf(\"follows a three-tier architecture\") !== -1\n);\nif (tierPara) {\n  tierPara.clear();\n  tierPara.insertText(\n    \"SwiftPay follows a three-tier architecture, which includes:\",\n    Word.InsertLocation.start\n  );\n}\n\n// --- 2) Append the new \"Prototype link\" section after the\n//        \"Microservices architecture...\" bullet paragraph.\n\nconst microPara = paragraphs.items.find(\n  (p) => p.text.indexOf(\"Microservices architecture can be adopted\") !== -1\n);\nif (microPara) {\n  const p1 = microPara.insertParagraph(\n    \"Prototype link\",\n    Word.InsertLocation.after\n  );\n\n  const p2 = p1.insertParagraph(\"               \", Word.InsertLocation.after);\n  p2.insertText(\"\u2022\", Word.InsertLocation.end);\n  p2.insertText(\n    \"        The interactive high-fidelity prototype for the swiftpay mobile application can be accessed via the following link:\",\n    Word.InsertLocation.end\n  );\n\n  const p3 = p2.insertParagraph(\"Figma link: [\", Word.InsertLocation.after);\n  p3.insertText(\n    \"https://www.figma.com/design/D7XU9E2VsMhGOzeUVeqh3I/Untitled?node-id=0-1&t=EyTmM9qg4n1mmEsj-1\",\n    Word.InsertLocation.end\n  );\n  p3.insertText(\"]\", Word.InsertLocation.end);\n\n  p3.insertParagraph(\n    \"This prototype demonstrates the min UI flow including login, registration, dashboard, transaction functions, and financial tips        \",\n    Word.InsertLocation.after\n  );\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document.\n\n$d = $word.ActiveDocument\n\n# --- 1) Collapse the \"SwiftPay\"-split runs (drops the spell-check\n#        proofErr markup by rewriting each paragraph as a single run).\n\n# Title paragraph: \"Phase 3: Software Architecture Design \u2013 SwiftPay\"\n$titleIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.IndexOf(\"Phase 3: Software Architecture Design\") -ge 0) {\n        $titleIdx = $i\n        break\n    }\n}\n$p = $d.Paragraphs.Item($titleIdx)\n$r = $p.Range\n$d.Range($r.Start, $r.End).Delete()\n$d.Range($r.Start, $r.Start).InsertBefore(\"Phase 3: Software Architecture Design \" + [char]8211 + \" SwiftPay\" + [char]13)\n\n# Overview paragraph.\n$overviewIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.IndexOf(\"The system architecture of\") -ge 0) {\n        $overviewIdx = $i\n        break\n    }\n}\n$p = $d.Paragraphs.Item($overviewIdx)\n$r = $p.Range\n$d.Range($r.Start, $r.End).Delete()\n$d.Range($r.Start, $r.Start).InsertBefore(\"The system architecture of SwiftPay outlines how different components of the mobile application will interact to deliver financial services such as registration, login, balance inquiry, money transfer, and more. The system is designed to be scalable, secure, and maintainable.\" + [char]13)\n\n# \"SwiftPay follows a three-tier architecture, which includes:\" paragraph.\n$tierIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.IndexOf(\"follows a three-tier architecture\") -ge 0) {\n        $tierIdx = $i\n        break\n    }\n}\n$p = $d.Paragraphs.Item($tierIdx)\n$r = $p.Range\n$d.Range($r.Start, $r.End).Delete()\n$d.Range($r.Start, $r.Start).InsertBefore(\"SwiftPay follows a three-tier architecture, which includes:\" + [char]13)\n\n# --- 2) Append the new \"Prototype link\" section after the\n#        \"Microservices architecture...\" bullet paragraph.\n\n$microIdx = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text.IndexOf(\"Microservices architecture can be adopted\") -ge 0) {\n        $microIdx = $i\n        break\n    }\n}\n$microPara = $d.Paragraphs.Item($microIdx)\n\n$microPara.Range.InsertParagraphAfter()\n$p1 = $microPara.Next()\n$p1.Range.Text = \"Prototype link\"\n\n$p1.Range.InsertParagraphAfter()\n$p2 = $p1.Next()\n$p2.Range.Text = \"               \" + [char]8226 + \"        The interactive high-fidelity prototype for the swiftpay mobile application can be accessed via the following link:\"\n\n$p2.Range.InsertParagraphAfter()\n$p3 = $p2.Next()\n$p3.Range.Text = \"Figma link: [https://www.figma.com/design/D7XU9E2VsMhGOzeUVeqh3I/Untitled?node-id=0-1&t=EyTmM9qg4n1mmEsj-1]\"\n\n$p3.Range.InsertParagraphAfter()\n$p4 = $p3.Next()\n$p4.Range.Text = \"This prototype demonstrates the min UI flow including login, registration, dashboard, transaction functions, and financial tips        \"\n"}
